$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "Neutrophils"
$ws.Range("G2").Value = 17.956883
$ws.Range("H2").Value = 35.913766
$ws.Range("I2").Value = 0.3392380274206944
$ws.Range("J2").Value = 0.2584869083704147
$ws.Range("M2").Value = 0.288097
$ws.Range("N2").Value = 0.864291
$ws.Range("Q2").Value = 5.173324121651
$ws.Range("R2").Value = 31.039944729906
$ws.Range("S2").Value = 0.3392380274206944
$ws.Range("T2").Value = 0.2584869083704147

# Row 3
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("I3").Value = 0.6187742881378531
$ws.Range("J3").Value = 0.7072248972319991
$ws.Range("M3").Value = 0.288097
$ws.Range("N3").Value = 0.864291
$ws.Range("Q3").Value = 9.436206120580998
$ws.Range("R3").Value = 84.925855085229
$ws.Range("S3").Value = 0.6187742881378531
$ws.Range("T3").Value = 0.7072248972319991

# Row 4
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("G4").Value = 1.903653
$ws.Range("H4").Value = 3.807306
$ws.Range("I4").Value = 0.03596345137480081
$ws.Range("J4").Value = 0.02740282812891664
$ws.Range("M4").Value = 0.288097
$ws.Range("N4").Value = 0.864291
$ws.Range("Q4").Value = 0.5484367183410001
$ws.Range("R4").Value = 3.290620310046
$ws.Range("S4").Value = 0.03596345137480081
$ws.Range("T4").Value = 0.02740282812891664

# Row 5
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3188806666666666
$ws.Range("H5").Value = 0.956642
$ws.Range("I5").Value = 0.006024233066651711
$ws.Range("J5").Value = 0.006885366268669519
$ws.Range("M5").Value = 0.288097
$ws.Range("N5").Value = 0.864291
$ws.Range("Q5").Value = 0.09186856342466666
$ws.Range("R5").Value = 0.8268170708220001
$ws.Range("S5").Value = 0.006024233066651711
$ws.Range("T5").Value = 0.006885366268669519
